# This script applies the edits described in the commit diff to DYNGROWTH_holdings.xlsx:
#  1. Updates the "as of" date in the confidential disclaimer text block (A81) from 2021-06-14 to 2021-07-07
#  2. Updates the Weight (column D) and Percent Change (column E) values for each holding row (2-78)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; unprotect so the cell values can be updated.
$ws.Unprotect()

# Update the confidential disclaimer date in the shared text block (A81)
$disclaimerText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-07-07 for illustrative purposes only and are subject to change."
$ws.Range("A81").Value = $disclaimerText

# Update Weight (D) and Percent Change (E) columns for holdings rows 2-78
$ws.Range("D2").Value = 0.0840448349190962
$ws.Range("E2").Value = 0
$ws.Range("D3").Value = 0.0501977448743165
$ws.Range("E3").Value = 0
$ws.Range("D4").Value = 0.04297445762588327
$ws.Range("E4").Value = 0
$ws.Range("D5").Value = 0.03599307252165512
$ws.Range("E5").Value = 0
$ws.Range("D6").Value = 0.03446627577407543
$ws.Range("E6").Value = 0
$ws.Range("D7").Value = 0.02982047089429782
$ws.Range("E7").Value = 0
$ws.Range("D8").Value = 0.02723556849374429
$ws.Range("E8").Value = 0
$ws.Range("D9").Value = 0.02795321759110845
$ws.Range("E9").Value = -0.006448839208942458
$ws.Range("D10").Value = 0.026617156632806
$ws.Range("E10").Value = 0
$ws.Range("D11").Value = 0.02690585505180274
$ws.Range("E11").Value = 0
$ws.Range("D12").Value = 0.02225344953015277
$ws.Range("E12").Value = 0
$ws.Range("D13").Value = 0.02213864388324174
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0.02025859899132041
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0.02054329150352564
$ws.Range("E15").Value = 0
$ws.Range("D16").Value = 0.02182855028024493
$ws.Range("E16").Value = 0
$ws.Range("D17").Value = 0.01891662021616688
$ws.Range("E17").Value = 0
$ws.Range("D18").Value = 0.01783684624918446
$ws.Range("E18").Value = 0
$ws.Range("D19").Value = 0.01712160980022882
$ws.Range("E19").Value = 0
$ws.Range("D20").Value = 0.01766504747383068
$ws.Range("E20").Value = 0
$ws.Range("D21").Value = 0.01508578976011972
$ws.Range("E21").Value = 0
$ws.Range("D22").Value = 0.01453324775290081
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 0.01347459583990995
$ws.Range("E23").Value = 0
$ws.Range("D24").Value = 0.01334458595585845
$ws.Range("E24").Value = 0
$ws.Range("D25").Value = 0.01196946740633884
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 0.01103991494403779
$ws.Range("E26").Value = 0
$ws.Range("D27").Value = 0.01181150721808298
$ws.Range("E27").Value = 0
$ws.Range("D28").Value = 0.01030137140102249
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0.01033451117538857
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0.01013722026590704
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0.01027150918675856
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 0.01021324145160943
$ws.Range("E32").Value = 0
$ws.Range("D33").Value = 0.01021123849821368
$ws.Range("E33").Value = 0
$ws.Range("D34").Value = 0.01022189056854563
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 0.01001585949879173
$ws.Range("E35").Value = 0
$ws.Range("D36").Value = 0.010932711415697
$ws.Range("E36").Value = 0
$ws.Range("D37").Value = 0.00815456953412153
$ws.Range("E37").Value = 0
$ws.Range("D38").Value = 0.009600383234177542
$ws.Range("E38").Value = 0
$ws.Range("D39").Value = 0.008492977614667368
$ws.Range("E39").Value = 0
$ws.Range("D40").Value = 0.009735992283403537
$ws.Range("E40").Value = 0
$ws.Range("D41").Value = 0.008520836875535549
$ws.Range("E41").Value = 0
$ws.Range("D42").Value = 0.008825240270021698
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 0.008955978500762572
$ws.Range("E43").Value = 0
$ws.Range("D44").Value = 0.008129987833355489
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 0.008998086043741438
$ws.Range("E45").Value = 0
$ws.Range("D46").Value = 0.008024195476725339
$ws.Range("E46").Value = 0
$ws.Range("D47").Value = 0.009414973480066268
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 0.008127256533270372
$ws.Range("E48").Value = 0
$ws.Range("D49").Value = 0.009007554550703173
$ws.Range("E49").Value = 0
$ws.Range("D50").Value = 0.006803850598695631
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0.00769461859978804
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 0.007902333971261092
$ws.Range("E52").Value = 0
$ws.Range("D53").Value = 0.008573368880505942
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 0.006608223506182619
$ws.Range("E54").Value = 0
$ws.Range("D55").Value = 0.006463348521417868
$ws.Range("E55").Value = 0
$ws.Range("D56").Value = 0.005589059364172323
$ws.Range("E56").Value = 0
$ws.Range("D57").Value = 0.006155257871816814
$ws.Range("E57").Value = 0
$ws.Range("D58").Value = 0.005905070784020213
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 0.006206879443425501
$ws.Range("E59").Value = 0
$ws.Range("D60").Value = 0.005147135010400596
$ws.Range("E60").Value = 0
$ws.Range("D61").Value = 0.005011890134519283
$ws.Range("E61").Value = 0
$ws.Range("D62").Value = 0.00531256075222244
$ws.Range("E62").Value = 0
$ws.Range("D63").Value = 0.004854430684612361
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0.004903229912799761
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0.004540331174824052
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 0.004281859143435934
$ws.Range("E66").Value = 0
$ws.Range("D67").Value = 0.003989336904320039
$ws.Range("E67").Value = 0
$ws.Range("D68").Value = 0.004023341590379729
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0.004129634685358816
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0.003996438284541339
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0.003350941031092321
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 0.00341753923150106
$ws.Range("E72").Value = 0
$ws.Range("D73").Value = 0.003768920987451199
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0.002717188368009301
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0.002465635630170142
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 0.001935035066968326
$ws.Range("E76").Value = 0
$ws.Range("D77").Value = 0.001594532989690564
$ws.Range("E77").Value = 0
$ws.Range("E78").Value = -0.0001802658056176965

# Restore sheet protection to its original (protected) state
$ws.Protect()
